{"js": "// The document contains a single table of two-digit-by-two-digit\n// multiplication problems (\"NN\u00d7NN=\") laid out five-per-row, with blank\n// spacer rows between each row of problems. This edit replaces the 25\n// problem strings, in document (reading) order, with a new set of\n// problems. Because a couple of the original strings repeat (e.g.\n// \"26\u00d769=\" appears twice), the replacement must be positional rather\n// than a simple text search/replace.\n\nconst replacements = [\n  \"21\u00d741=\", \"73\u00d726=\", \"70\u00d772=\", \"42\u00d781=\", \"91\u00d766=\",\n  \"54\u00d783=\", \"79\u00d711=\", \"52\u00d778=\", \"83\u00d784=\", \"29\u00d768=\",\n  \"26\u00d738=\", \"51\u00d767=\", \"62\u00d796=\", \"88\u00d799=\", \"36\u00d761=\",\n  \"15\u00d745=\", \"79\u00d716=\", \"22\u00d794=\", \"78\u00d725=\", \"11\u00d720=\",\n  \"33\u00d797=\", \"92\u00d764=\", \"13\u00d742=\", \"74\u00d711=\", \"32\u00d799=\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst rows = table.values;\nlet i = 0;\nfor (let r = 0; r < rows.length; r++) {\n  for (let c = 0; c < rows[r].length; c++) {\n    const text = rows[r][c];\n    // Only the cells holding a \"NN\u00d7NN=\" style expression are replaced;\n    // the interleaved blank spacer rows are left untouched.\n    if (text && text.indexOf(\"\u00d7\") !== -1) {\n      if (i < replacements.length) {\n        table.getCell(r, c).value = replacements[i];\n        i++;\n      }\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table of two-digit-by-two-digit\n# multiplication problems (\"NN\u00d7NN=\") laid out five-per-row, with blank\n# spacer rows between each row of problems. This edit replaces the 25\n# problem strings, in document (reading) order, with a new set of\n# problems. Because a couple of the original strings repeat (e.g.\n# \"26\u00d769=\" appears twice), the replacement must be positional rather\n# than a simple text search/replace.\n\n$replacements = @(\n    \"21\u00d741=\", \"73\u00d726=\", \"70\u00d772=\", \"42\u00d781=\", \"91\u00d766=\",\n    \"54\u00d783=\", \"79\u00d711=\", \"52\u00d778=\", \"83\u00d784=\", \"29\u00d768=\",\n    \"26\u00d738=\", \"51\u00d767=\", \"62\u00d796=\", \"88\u00d799=\", \"36\u00d761=\",\n    \"15\u00d745=\", \"79\u00d716=\", \"22\u00d794=\", \"78\u00d725=\", \"11\u00d720=\",\n    \"33\u00d797=\", \"92\u00d764=\", \"13\u00d742=\", \"74\u00d711=\", \"32\u00d799=\"\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $text = $cell.Range.Text\n        # Only the cells holding a \"NN\u00d7NN=\" style expression are replaced;\n        # the interleaved blank spacer rows are left untouched.\n        if ($text -like \"*\u00d7*\") {\n            if ($i -lt $replacements.Count) {\n                $cell.Range.Text = $replacements[$i]\n                $i++\n            }\n        }\n    }\n}\n"}
